$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old"/"_new" suffixes become "_FV2410"/"_FV2504"
$headers = @(
  "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
  "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410",
  "diff",
  "Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504",
  "Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Freeze the header row (row 1) - pane split with topLeftCell A2
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the used range into an Excel Table ("Table1") spanning A1:U55
$rng = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Type]::Missing, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Output "done"
